# Added Cross Browser Testing
# - Insert a new "Browser" column into the TestData sheet (between
#   "Execution Flag" and "userName") with chrome/firefox values per row.
# - Fix the stray "admin124" typo back to "admin123".
# - Make TestData the active/selected sheet with cell E9 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Insert new column C ("Browser"), shifting userName/Password/fname right.
$ws.Columns("C").Insert()
# Keep the new column's width in line with its neighbours (matches the
# author's sheet, where the inserted column picked up the existing width).
$ws.Columns("C").ColumnWidth = 8.6

$ws.Range("C1").Value = "Browser"
$ws.Range("C2").Value = "chrome"
$ws.Range("C3").Value = "firefox"
$ws.Range("C4").Value = "chrome"
$ws.Range("C5").Value = "firefox"
$ws.Range("C6").Value = "chrome"

# Correct the old "admin124" typo (old D3, now E3) to "admin123".
$ws.Range("E3").Value = "admin123"

# Make TestData the active sheet / tab, with E9 selected.
$ws.Activate()
$ws.Range("E9").Select()
